$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" data table (B15:J29) gets sorted by column E (Periodo
# Mora) in ascending order. Re-apply the same net effect by writing the
# sorted "Periodo Mora" labels and their matching "Valor Mora" amounts back
# into E16:E29 / F16:F29 (the only columns whose contents actually change
# when the table is re-sorted ascending instead of descending).

$periodos = @("2310","2311","2312","2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411")
$valores  = @(46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,46400,21654)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
